$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2026-02-09 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-02-10 Tuesday", 2)

# Update the table cells using the Tables/Cell object model so that
# duplicate text values (e.g. "549÷3=183, 0") are each replaced with the
# correct, position-specific new value.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "972÷9=108, 0"
$t.Cell(1,2).Range.Text = "227÷7=32, 3"
$t.Cell(1,3).Range.Text = "782÷9=86, 8"
$t.Cell(1,4).Range.Text = "357÷7=51, 0"
$t.Cell(1,5).Range.Text = "968÷7=138, 2"

$t.Cell(5,1).Range.Text = "384÷5=76, 4"
$t.Cell(5,2).Range.Text = "294÷8=36, 6"
$t.Cell(5,3).Range.Text = "895÷3=298, 1"
$t.Cell(5,4).Range.Text = "778÷4=194, 2"
$t.Cell(5,5).Range.Text = "845÷8=105, 5"

$t.Cell(9,1).Range.Text = "477÷3=159, 0"
$t.Cell(9,2).Range.Text = "950÷5=190, 0"
$t.Cell(9,3).Range.Text = "938÷6=156, 2"
$t.Cell(9,4).Range.Text = "964÷3=321, 1"
$t.Cell(9,5).Range.Text = "852÷9=94, 6"

$t.Cell(13,1).Range.Text = "441÷9=49, 0"
$t.Cell(13,2).Range.Text = "660÷6=110, 0"
$t.Cell(13,3).Range.Text = "380÷6=63, 2"
$t.Cell(13,4).Range.Text = "144÷7=20, 4"
$t.Cell(13,5).Range.Text = "116÷5=23, 1"

$t.Cell(17,1).Range.Text = "248÷8=31, 0"
$t.Cell(17,2).Range.Text = "308÷3=102, 2"
$t.Cell(17,3).Range.Text = "554÷9=61, 5"
$t.Cell(17,4).Range.Text = "703÷6=117, 1"
$t.Cell(17,5).Range.Text = "147÷3=49, 0"
